# Insert a new ".NET Native (NGEN with a Turbo)" slide right before the
# existing "Platforms, Frameworks and Devices oh my..." slide (currently
# slide 13), pushing it and everything after it down by one position.

$p = $ppt.ActivePresentation

# Reuse the "Title and Content" custom layout already used by the
# neighboring slide so the new slide matches the deck's look & feel.
$refSlide = $p.Slides.Item(13)
$layout = $refSlide.CustomLayout

$s = $p.Slides.AddSlide(13, $layout)

# Title
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = ".NET Native (NGEN with a Turbo)"
$title.LanguageID = "en-US"

# Body / content placeholder
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Only on Universal Windows Projects`rTypical path:`rC# -> JIT -> IL`r.NET Native path:`rC# -> IL -> Native`rCommon Issues:`rReflection`rSerialization/Deserialization`r"
$body.LanguageID = "en-US"

$body.Paragraphs(3).IndentLevel = 2
$body.Paragraphs(5).IndentLevel = 2
$body.Paragraphs(7).IndentLevel = 2
$body.Paragraphs(8).IndentLevel = 2
$body.Paragraphs(9).IndentLevel = 2
